$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-22 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-23 Saturday", 2) | Out-Null
$d.Content.Find.Execute("38×88=3344", $true, $false, $false, $false, $false, $true, 1, $false, "43×77=3311", 2) | Out-Null
$d.Content.Find.Execute("68×25=1700", $true, $false, $false, $false, $false, $true, 1, $false, "77×99=7623", 2) | Out-Null
$d.Content.Find.Execute("59×54=3186", $true, $false, $false, $false, $false, $true, 1, $false, "48×55=2640", 2) | Out-Null
$d.Content.Find.Execute("77×61=4697", $true, $false, $false, $false, $false, $true, 1, $false, "50×97=4850", 2) | Out-Null
$d.Content.Find.Execute("83×19=1577", $true, $false, $false, $false, $false, $true, 1, $false, "99×30=2970", 2) | Out-Null
$d.Content.Find.Execute("20×95=1900", $true, $false, $false, $false, $false, $true, 1, $false, "54×14=756", 2) | Out-Null
$d.Content.Find.Execute("34×60=2040", $true, $false, $false, $false, $false, $true, 1, $false, "62×96=5952", 2) | Out-Null
$d.Content.Find.Execute("63×78=4914", $true, $false, $false, $false, $false, $true, 1, $false, "25×94=2350", 2) | Out-Null
$d.Content.Find.Execute("87×36=3132", $true, $false, $false, $false, $false, $true, 1, $false, "46×38=1748", 2) | Out-Null
$d.Content.Find.Execute("71×40=2840", $true, $false, $false, $false, $false, $true, 1, $false, "42×69=2898", 2) | Out-Null
$d.Content.Find.Execute("27×49=1323", $true, $false, $false, $false, $false, $true, 1, $false, "29×55=1595", 2) | Out-Null
$d.Content.Find.Execute("16×16=256", $true, $false, $false, $false, $false, $true, 1, $false, "36×30=1080", 2) | Out-Null
$d.Content.Find.Execute("66×84=5544", $true, $false, $false, $false, $false, $true, 1, $false, "65×91=5915", 2) | Out-Null
$d.Content.Find.Execute("59×12=708", $true, $false, $false, $false, $false, $true, 1, $false, "15×58=870", 2) | Out-Null
$d.Content.Find.Execute("71×93=6603", $true, $false, $false, $false, $false, $true, 1, $false, "39×43=1677", 2) | Out-Null
$d.Content.Find.Execute("27×51=1377", $true, $false, $false, $false, $false, $true, 1, $false, "44×26=1144", 2) | Out-Null
$d.Content.Find.Execute("59×90=5310", $true, $false, $false, $false, $false, $true, 1, $false, "55×46=2530", 2) | Out-Null
$d.Content.Find.Execute("54×40=2160", $true, $false, $false, $false, $false, $true, 1, $false, "23×48=1104", 2) | Out-Null
$d.Content.Find.Execute("49×72=3528", $true, $false, $false, $false, $false, $true, 1, $false, "17×34=578", 2) | Out-Null
$d.Content.Find.Execute("94×33=3102", $true, $false, $false, $false, $false, $true, 1, $false, "68×37=2516", 2) | Out-Null
$d.Content.Find.Execute("89×59=5251", $true, $false, $false, $false, $false, $true, 1, $false, "80×68=5440", 2) | Out-Null
$d.Content.Find.Execute("25×35=875", $true, $false, $false, $false, $false, $true, 1, $false, "60×11=660", 2) | Out-Null
$d.Content.Find.Execute("69×70=4830", $true, $false, $false, $false, $false, $true, 1, $false, "17×99=1683", 2) | Out-Null
$d.Content.Find.Execute("42×82=3444", $true, $false, $false, $false, $false, $true, 1, $false, "71×77=5467", 2) | Out-Null
$d.Content.Find.Execute("85×94=7990", $true, $false, $false, $false, $false, $true, 1, $false, "58×15=870", 2) | Out-Null
